$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell without triggering Excel's
# automatic date-recognition (which would otherwise store the date
# strings "01-0X-2021" as serial date numbers with a new number format).
# We build the literal text with a formula on a scratch cell (far outside
# the used range), copy it, and paste-special only the resulting VALUE
# into the destination cell; this keeps the destination a plain shared
# string and leaves styles.xml untouched.
function Set-TextValue($row, $col, $text) {
    $scratch = $ws.Cells.Item(1000, 1000)
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)
    $scratch.Clear()
    $excel.CutCopyMode = $false
}

# --- Update existing rows 98-100 with revised monthly figures ---
# Row 98
$ws.Cells.Item(98,2).Value = 633
$ws.Cells.Item(98,3).Value = 57
$ws.Cells.Item(98,4).Value = 1208
$ws.Cells.Item(98,5).Value = 1080
$ws.Cells.Item(98,6).Value = 128
$ws.Cells.Item(98,7).Value = 1332
$ws.Cells.Item(98,8).Value = 369
$ws.Cells.Item(98,9).Value = 174
$ws.Cells.Item(98,10).Value = 29
$ws.Cells.Item(98,11).Value = 69
$ws.Cells.Item(98,12).Value = 94
$ws.Cells.Item(98,13).Value = 101
$ws.Cells.Item(98,14).Value = 197
$ws.Cells.Item(98,15).Value = 76
$ws.Cells.Item(98,16).Value = 218
$ws.Cells.Item(98,17).Value = 371
$ws.Cells.Item(98,18).Value = 634
$ws.Cells.Item(98,19).Value = 1351
$ws.Cells.Item(98,20).Value = 185
$ws.Cells.Item(98,21).Value = 567
$ws.Cells.Item(98,22).Value = 500
$ws.Cells.Item(98,23).Value = 740
$ws.Cells.Item(98,24).Value = 1209
$ws.Cells.Item(98,25).Value = 1001
$ws.Cells.Item(98,26).Value = 859
$ws.Cells.Item(98,27).Value = 626
$ws.Cells.Item(98,28).Value = 11173
$ws.Cells.Item(98,29).Value = 1040
$ws.Cells.Item(98,30).Value = 52
$ws.Cells.Item(98,31).Value = 12264

# Row 99
$ws.Cells.Item(99,2).Value = 701
$ws.Cells.Item(99,3).Value = 59
$ws.Cells.Item(99,4).Value = 1112
$ws.Cells.Item(99,5).Value = 1001
$ws.Cells.Item(99,6).Value = 111
$ws.Cells.Item(99,7).Value = 1299
$ws.Cells.Item(99,8).Value = 353
$ws.Cells.Item(99,9).Value = 168
$ws.Cells.Item(99,10).Value = 27
$ws.Cells.Item(99,11).Value = 68
$ws.Cells.Item(99,12).Value = 85
$ws.Cells.Item(99,13).Value = 104
$ws.Cells.Item(99,14).Value = 208
$ws.Cells.Item(99,15).Value = 74
$ws.Cells.Item(99,16).Value = 210
$ws.Cells.Item(99,17).Value = 331
$ws.Cells.Item(99,18).Value = 633
$ws.Cells.Item(99,19).Value = 1319
$ws.Cells.Item(99,20).Value = 197
$ws.Cells.Item(99,21).Value = 574
$ws.Cells.Item(99,22).Value = 461
$ws.Cells.Item(99,23).Value = 740
$ws.Cells.Item(99,24).Value = 1170
$ws.Cells.Item(99,25).Value = 993
$ws.Cells.Item(99,26).Value = 772
$ws.Cells.Item(99,27).Value = 626
$ws.Cells.Item(99,28).Value = 10885
$ws.Cells.Item(99,29).Value = 1018
$ws.Cells.Item(99,30).Value = 67
$ws.Cells.Item(99,31).Value = 11968

# Row 100
$ws.Cells.Item(100,2).Value = 639
$ws.Cells.Item(100,3).Value = 113
$ws.Cells.Item(100,4).Value = 1280
$ws.Cells.Item(100,5).Value = 1155
$ws.Cells.Item(100,6).Value = 124
$ws.Cells.Item(100,7).Value = 1456
$ws.Cells.Item(100,8).Value = 412
$ws.Cells.Item(100,9).Value = 177
$ws.Cells.Item(100,10).Value = 32
$ws.Cells.Item(100,11).Value = 79
$ws.Cells.Item(100,12).Value = 99
$ws.Cells.Item(100,13).Value = 106
$ws.Cells.Item(100,14).Value = 220
$ws.Cells.Item(100,15).Value = 81
$ws.Cells.Item(100,16).Value = 247
$ws.Cells.Item(100,17).Value = 353
$ws.Cells.Item(100,18).Value = 732
$ws.Cells.Item(100,19).Value = 1491
$ws.Cells.Item(100,20).Value = 196
$ws.Cells.Item(100,21).Value = 572
$ws.Cells.Item(100,22).Value = 469
$ws.Cells.Item(100,23).Value = 734
$ws.Cells.Item(100,24).Value = 1249
$ws.Cells.Item(100,25).Value = 969
$ws.Cells.Item(100,26).Value = 1683
$ws.Cells.Item(100,27).Value = 630
$ws.Cells.Item(100,28).Value = 12603
$ws.Cells.Item(100,29).Value = 1151
$ws.Cells.Item(100,30).Value = 87
$ws.Cells.Item(100,31).Value = 13843

# --- Append new rows 101-103 (Apr/May/Jun 2021) ---
# Row 101
Set-TextValue 101 1 "01-04-2021"
$ws.Cells.Item(101,2).Value = 378
$ws.Cells.Item(101,3).Value = 93
$ws.Cells.Item(101,4).Value = 1240
$ws.Cells.Item(101,5).Value = 1111
$ws.Cells.Item(101,6).Value = 128
$ws.Cells.Item(101,7).Value = 1373
$ws.Cells.Item(101,8).Value = 380
$ws.Cells.Item(101,9).Value = 193
$ws.Cells.Item(101,10).Value = 28
$ws.Cells.Item(101,11).Value = 72
$ws.Cells.Item(101,12).Value = 92
$ws.Cells.Item(101,13).Value = 79
$ws.Cells.Item(101,14).Value = 211
$ws.Cells.Item(101,15).Value = 84
$ws.Cells.Item(101,16).Value = 237
$ws.Cells.Item(101,17).Value = 331
$ws.Cells.Item(101,18).Value = 699
$ws.Cells.Item(101,19).Value = 1297
$ws.Cells.Item(101,20).Value = 177
$ws.Cells.Item(101,21).Value = 504
$ws.Cells.Item(101,22).Value = 469
$ws.Cells.Item(101,23).Value = 729
$ws.Cells.Item(101,24).Value = 1298
$ws.Cells.Item(101,25).Value = 962
$ws.Cells.Item(101,26).Value = 1678
$ws.Cells.Item(101,27).Value = 632
$ws.Cells.Item(101,28).Value = 11863
$ws.Cells.Item(101,29).Value = 1040
$ws.Cells.Item(101,30).Value = 59
$ws.Cells.Item(101,31).Value = 12965

# Row 102
Set-TextValue 102 1 "01-05-2021"
$ws.Cells.Item(102,2).Value = 279
$ws.Cells.Item(102,3).Value = 102
$ws.Cells.Item(102,4).Value = 1271
$ws.Cells.Item(102,5).Value = 1133
$ws.Cells.Item(102,6).Value = 138
$ws.Cells.Item(102,7).Value = 1387
$ws.Cells.Item(102,8).Value = 358
$ws.Cells.Item(102,9).Value = 206
$ws.Cells.Item(102,10).Value = 32
$ws.Cells.Item(102,11).Value = 75
$ws.Cells.Item(102,12).Value = 98
$ws.Cells.Item(102,13).Value = 93
$ws.Cells.Item(102,14).Value = 198
$ws.Cells.Item(102,15).Value = 82
$ws.Cells.Item(102,16).Value = 243
$ws.Cells.Item(102,17).Value = 365
$ws.Cells.Item(102,18).Value = 736
$ws.Cells.Item(102,19).Value = 1404
$ws.Cells.Item(102,20).Value = 212
$ws.Cells.Item(102,21).Value = 523
$ws.Cells.Item(102,22).Value = 489
$ws.Cells.Item(102,23).Value = 746
$ws.Cells.Item(102,24).Value = 1291
$ws.Cells.Item(102,25).Value = 963
$ws.Cells.Item(102,26).Value = 1733
$ws.Cells.Item(102,27).Value = 636
$ws.Cells.Item(102,28).Value = 12128
$ws.Cells.Item(102,29).Value = 1219
$ws.Cells.Item(102,30).Value = 51
$ws.Cells.Item(102,31).Value = 13392

# Row 103
Set-TextValue 103 1 "01-06-2021"
$ws.Cells.Item(103,2).Value = 213
$ws.Cells.Item(103,3).Value = 76
$ws.Cells.Item(103,4).Value = 1264
$ws.Cells.Item(103,5).Value = 1124
$ws.Cells.Item(103,6).Value = 140
$ws.Cells.Item(103,7).Value = 1385
$ws.Cells.Item(103,8).Value = 352
$ws.Cells.Item(103,9).Value = 201
$ws.Cells.Item(103,10).Value = 29
$ws.Cells.Item(103,11).Value = 73
$ws.Cells.Item(103,12).Value = 97
$ws.Cells.Item(103,13).Value = 98
$ws.Cells.Item(103,14).Value = 206
$ws.Cells.Item(103,15).Value = 82
$ws.Cells.Item(103,16).Value = 246
$ws.Cells.Item(103,17).Value = 378
$ws.Cells.Item(103,18).Value = 725
$ws.Cells.Item(103,19).Value = 1365
$ws.Cells.Item(103,20).Value = 222
$ws.Cells.Item(103,21).Value = 538
$ws.Cells.Item(103,22).Value = 469
$ws.Cells.Item(103,23).Value = 749
$ws.Cells.Item(103,24).Value = 1313
$ws.Cells.Item(103,25).Value = 977
$ws.Cells.Item(103,26).Value = 1765
$ws.Cells.Item(103,27).Value = 636
$ws.Cells.Item(103,28).Value = 12075
$ws.Cells.Item(103,29).Value = 1195
$ws.Cells.Item(103,30).Value = 55
$ws.Cells.Item(103,31).Value = 13320

